$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Heading text: "Data Clumps" -> "Long Methods".
#    A plain Find/Replace regenerates the run from scratch and silently
#    drops the <w:lastRenderedPageBreak/> marker that sits in front of
#    the text in that run. To keep it, round-trip the paragraph's OOXML
#    (Range.WordOpenXML) - which preserves every attribute/property
#    exactly as authored - patch just the <w:t> run, and push it back
#    with InsertXML (which replaces only that range's contents).
# ---------------------------------------------------------------------
$didHeading = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Data Clumps*") {
        $xml = $p.Range.WordOpenXML
        $old = "<w:r><w:t>Data Clumps</w:t></w:r>"
        $new = "<w:r><w:lastRenderedPageBreak/><w:t>Long Methods</w:t></w:r>"
        if ($xml.Contains($old)) {
            $xml = $xml.Replace($old, $new)
            $p.Range.InsertXML($xml)
            $didHeading = $true
        }
    }
}
if (-not $didHeading) {
    # Fallback: plain text replace if the expected run shape wasn't found.
    $d.Content.Find.Execute("Data Clumps", $true, $false, $false, $false, $false,
                             $true, 1, $false, "Long Methods", 2)
}

# ---------------------------------------------------------------------
# 2. Collapse the redundant run splits in the two "Folder Encontrada:
#    ...\o" + "rg\..." paragraphs (imgscalr\Scalr.java and
#    ganttProject\WebStartIDClass.java) into their natural merged form,
#    i.e. the first run keeps the trailing "o" and the remainder of the
#    path becomes a single sibling run. Word's WordOpenXML projection
#    already represents the paragraph with adjacent equally-formatted
#    runs coalesced, so round-tripping Range.WordOpenXML back through
#    InsertXML merges the runs while leaving the visible text (and the
#    surviving runs' rsidRPr/rPr) exactly as authored.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if (($t -like "Folder Encontrada: ganttproject\src\main\java\org\imgscalr\Scalr.java*") -or
        ($t -like "Folder Encontrada: ganttproject\src\main\java\org\ganttProject\WebStartIDClass.java*")) {
        $p.Range.InsertXML($p.Range.WordOpenXML)
    }
}
